$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder the attribute rows (2-6) so that the Attribute/Type pairs end up
# in the new order described by the commit (sensor row count attribute
# removed from query output, shuffling remaining attribute rows).
$ws.Range("A2").Value = "concept:name"
$ws.Range("B2").Value = "str"

$ws.Range("A3").Value = "org:resource"
$ws.Range("B3").Value = "str"

$ws.Range("A4").Value = "time:timestamp"
$ws.Range("B4").Value = "datetime"

$ws.Range("A5").Value = "SubProcessID"
$ws.Range("B5").Value = "str"

$ws.Range("A6").Value = "operation_end_time"
$ws.Range("B6").Value = "datetime"
